$wb = $excel.ActiveWorkbook

# --- Rushing sheet: add D.Samuel as row 12 ---
$wsRushing = $wb.Worksheets.Item("Rushing")

# Copy formatting from the row above (A11 carries the "row number" style)
# then override the borders so only left/right remain (matches the new
# style used for the last "total-like" row of new data).
$wsRushing.Range("A11").Copy()
$wsRushing.Range("A12").PasteSpecial(-4122)
$wsRushing.Range("A12").Borders.Item(8).LineStyle = -4142
$wsRushing.Range("A12").Borders.Item(9).LineStyle = -4142

$wsRushing.Range("A12").Value = 10
$wsRushing.Range("B12").Value = "D.Samuel"
$wsRushing.Range("C12").Value = 12
$wsRushing.Range("D12").Value = 11
$wsRushing.Range("E12").Value = 3
$wsRushing.Range("F12").Value = 7

# --- Receiving sheet: add D.Samuel as row 14 ---
$wsReceiving = $wb.Worksheets.Item("Receiving")

$wsReceiving.Range("A13").Copy()
$wsReceiving.Range("A14").PasteSpecial(-4122)
$wsReceiving.Range("A14").Borders.Item(8).LineStyle = -4142
$wsReceiving.Range("A14").Borders.Item(9).LineStyle = -4142

$wsReceiving.Range("A14").Value = 12
$wsReceiving.Range("B14").Value = "D.Samuel"
$wsReceiving.Range("C14").Value = 69
$wsReceiving.Range("D14").Value = 40
$wsReceiving.Range("E14").Value = 22
$wsReceiving.Range("F14").Value = 15
$wsReceiving.Range("G14").Value = 8
$wsReceiving.Range("H14").Value = 4

# Make Receiving the active sheet (tabSelected in final file)
$wsReceiving.Activate()
